# Apply updates to GAIA_TESS_candidate_matches_unique_planets worksheet
# - rows 8-12 are cyclically reordered (row that was 9 moves to 8, 10->9, 11->10, 12->11, 8->12)
# - column Q (HZ_Detection_Limit) values are recalculated/updated across rows 2-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("Q2").Value = 0.4573291322971118

# Row 3
$ws.Range("Q3").Value = 0.4573291322971118

# Row 4
$ws.Range("Q4").Value = 0.4573291322971118

# Row 5
$ws.Range("Q5").Value = 0.5747214936203638

# Row 6
$ws.Range("Q6").Value = 0.5747214936203638

# Row 7
$ws.Range("Q7").Value = 0.6689697604315441

# Row 8
$ws.Range("A8").Value = 5789.01
$ws.Range("B8").Value = 87216634
$ws.Range("C8").Value = 302.773447
$ws.Range("D8").Value = 16.189714
$ws.Range("E8").Value = 2459776.255982
$ws.Range("F8").Value = 12.9256051
$ws.Range("G8").Value = 1.691
$ws.Range("H8").Value = 790
$ws.Range("I8").Value = 2.53704
$ws.Range("J8").Value = 8.719390000000001
$ws.Range("K8").Value = 478
$ws.Range("L8").Value = 6.5715
$ws.Range("M8").Value = 20.4581
$ws.Range("N8").Value = 5132
$ws.Range("O8").Value = 4.49
$ws.Range("P8").Value = 0.88
$ws.Range("Q8").Value = 0.8326105656603114
$ws.Range("R8").Value = 0.2860097255552455

# Row 9
$ws.Range("A9").Value = 2011.02
$ws.Range("B9").Value = 136916387
$ws.Range("C9").Value = 230.440115
$ws.Range("D9").Value = -48.318817
$ws.Range("E9").Value = 2458650.89757
$ws.Range("F9").Value = 27.5920718
$ws.Range("G9").Value = 3.2743016
$ws.Range("H9").Value = 555.3352929
$ws.Range("I9").Value = 2.7050278
$ws.Range("J9").Value = 31.3830711
$ws.Range("K9").Value = 603.6616824
$ws.Range("L9").Value = 5.0494
$ws.Range("M9").Value = 14.682
$ws.Range("N9").Value = 5739
$ws.Range("O9").Value = 4.51
$ws.Range("P9").Value = 1.02607
$ws.Range("Q9").Value = 0.8786003547184101
$ws.Range("R9").Value = 0.8241668035706786

# Row 10
$ws.Range("A10").Value = 2011.03
$ws.Range("E10").Value = 2460081.126376
$ws.Range("G10").Value = 3.5341354
$ws.Range("H10").Value = 195.1972939
$ws.Range("I10").Value = 1.6208326
$ws.Range("J10").Value = 335.6075222
$ws.Range("K10").Value = 1091.635552
$ws.Range("Q10").Value = 0.8786003547184101

# Row 11
$ws.Range("A11").Value = 2011.01
$ws.Range("E11").Value = 2458631.767617
$ws.Range("F11").Value = 11.5778362
$ws.Range("G11").Value = 3.8772499
$ws.Range("H11").Value = 217.3994918
$ws.Range("I11").Value = 1.5422129
$ws.Range("J11").Value = 99.9013415
$ws.Range("K11").Value = 806.3295568999999
$ws.Range("Q11").Value = 0.8786003547184101

# Row 12
$ws.Range("A12").Value = 1099.01
$ws.Range("B12").Value = 290348383
$ws.Range("C12").Value = 328.718171
$ws.Range("D12").Value = -77.338802
$ws.Range("E12").Value = 2460140.220927
$ws.Range("F12").Value = 6.441006
$ws.Range("G12").Value = 1.812
$ws.Range("H12").Value = 914
$ws.Range("I12").Value = 2.56374
$ws.Range("J12").Value = 87.89019999999999
$ws.Range("K12").Value = 853
$ws.Range("L12").Value = 7.3661
$ws.Range("M12").Value = 23.606
$ws.Range("N12").Value = 4867
$ws.Range("O12").Value = 4.438
$ws.Range("P12").Value = 0.8
$ws.Range("Q12").Value = 0.9099502856655908
$ws.Range("R12").Value = 0.1472633162803166

# Row 13
$ws.Range("Q13").Value = 1.386957062992251

# Row 14
$ws.Range("Q14").Value = 1.546870014982985

# Row 15
$ws.Range("Q15").Value = 1.546870014982985

# Row 16
$ws.Range("Q16").Value = 1.692375571059443

# Row 17
$ws.Range("Q17").Value = 4.872880388587285

# Row 18
$ws.Range("Q18").Value = 4.872880388587285

# Row 19
$ws.Range("Q19").Value = 4.872880388587285

# Row 20
$ws.Range("Q20").Value = 10.16437674197399

# Row 21
$ws.Range("Q21").Value = 10.16437674197399

# Clear cells that no longer have a value after the row shift
$ws.Range("F10").ClearContents()

Write-Host "Applied GAIA/TESS candidate matches update"
